# Update res_bus vm_pu results for the 380 kV case (rows 2-25).
# Values taken from the new power-flow run; column layout unchanged
# (A = bus index, B:F and I:N = per-case voltage magnitudes, G stays 1, H stays blank).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$bfData = @(
    @(1.02,1.023641669795336,1.035515722885704,1.024234739284227,1.046515598926637),
    @(1.02,1.024509501867758,1.036204942319545,1.024968729772146,1.047366565598113),
    @(1.02,1.025071496332558,1.036650901360054,1.025444460184685,1.04791747081863),
    @(1.02,1.02530786492205,1.036838377918205,1.025644644596885,1.048149134855375),
    @(1.02,1.025347558430472,1.03686985573534,1.025678267378374,1.048188035935351),
    @(1.02,1.025074654286002,1.036653406451367,1.025447134324544,1.047920566077716),
    @(1.02,1.023934863882948,1.035748649261907,1.024482630342237,1.04680312948708),
    @(1.02,1.021929908546331,1.034154336249661,1.02278917071797,1.044836234691915),
    @(1.02,1.020595707599092,1.033091553631536,1.021664410349616,1.043526540149133),
    @(1.02,1.020018578661994,1.032631400560356,1.021178396518055,1.042959822731072),
    @(1.02,1.019804297089977,1.03246048677105,1.020998023468486,1.042749378892454),
    @(1.02,1.019850257128877,1.032497147989877,1.021036707098711,1.042794517055105),
    @(1.02,1.020000864240519,1.032617272609291,1.021163483671774,1.042942426138344),
    @(1.02,1.020093670193698,1.03269128640361,1.021241615375542,1.043033565842946),
    @(1.02,1.02063402226723,1.033122093439982,1.021696687010108,1.043564159667016),
    @(1.02,1.020973129451843,1.033392338918958,1.021982414296028,1.043897092567674),
    @(1.02,1.02117098166366,1.033549972164942,1.022149171993067,1.044091324104147),
    @(1.02,1.021238453736278,1.033603721591615,1.022206048597245,1.044157558328754),
    @(1.02,1.020936740582841,1.033363343748713,1.021951748333884,1.043861368135888),
    @(1.02,1.019956511709915,1.032581898680727,1.02112614686285,1.042898868922368),
    @(1.02,1.01934072319999,1.032090617677064,1.020607951001597,1.042294056906632),
    @(1.02,1.019667114575401,1.032351050239802,1.020882571275296,1.042614645474683),
    @(1.02,1.02095318297125,1.03337644540704,1.021965604660172,1.043877510353479),
    @(1.02,1.022447814217772,1.034566495551391,1.023226235278996,1.045344455604575)
)

$inData = @(
    @(1.034887709612122,1.028820747242331,1.038312105317031,1.027063837507516,1.049280817864134,1.013661601214851),
    @(1.035062272103598,1.029327568026764,1.038811055226286,1.027605074918201,1.049943340577621,1.01382978519807),
    @(1.03517355846195,1.029655336199618,1.039133221108823,1.027955415151753,1.050371666395979,1.01393852491798),
    @(1.035219943246847,1.029793085962322,1.039268493550763,1.028102726473607,1.050551644209011,1.013984217948532),
    @(1.03522770797512,1.029816212140389,1.039291196609197,1.028127462314756,1.050581857929753,1.013991888756953),
    @(1.035174179830526,1.029657176992977,1.039135029280182,1.027957383423286,1.050374071625427,1.013939135554121),
    @(1.034947048523008,1.028992066004948,1.038480869423264,1.027246724907318,1.049504796296393,1.013718457448211),
    @(1.034534082726626,1.027818738336733,1.03732294659817,1.02599545486541,1.04797025442903,1.013328950146036),
    @(1.034250262658162,1.027035703936103,1.036547576028559,1.025162022566047,1.04694546107253,1.013068870475084),
    @(1.034125357537085,1.026696462386958,1.036211037162423,1.024801330158779,1.046501315385941,1.012956161536344),
    @(1.034078661166789,1.02657042670485,1.036085913035204,1.024667382622744,1.04633628087008,1.012914282902814),
    @(1.034088691306822,1.026597462956649,1.036112757931244,1.024696113460672,1.046371683980067,1.012923266616802),
    @(1.034121503737161,1.026686044767332,1.036200696774945,1.024790257400922,1.046487674776531,1.012952700108741),
    @(1.034141680690873,1.026740619498328,1.036254863121796,1.024848266560297,1.046559132759065,1.012970833276054),
    @(1.034258509937923,1.027058214533941,1.036569894298962,1.025185964627587,1.046974929194278,1.013076348678038),
    @(1.034331256499462,1.02725738546528,1.036767292245621,1.025397845259719,1.047235640423889,1.013142511170633),
    @(1.03437349442946,1.027373540815531,1.036882354015813,1.025521449826854,1.047387669845641,1.013181093658764),
    @(1.034387863557227,1.027413143753637,1.036921573985962,1.025563598835536,1.047439501281946,1.01319424775528),
    @(1.034323471532906,1.027236018130125,1.036746121280592,1.025375110583736,1.047207672601709,1.013135413486593),
    @(1.034111849598162,1.026659960347775,1.036174804240822,1.024762533527577,1.04645351999493,1.012944033048733),
    @(1.033977053120999,1.026297618556854,1.03581490896691,1.02437755437599,1.045979013552052,1.012823626581052),
    @(1.03404867607698,1.026489716763402,1.036005760803909,1.02458162235528,1.04623059017178,1.012887463581862),
    @(1.034326989823415,1.0272456731639,1.036755687766798,1.025385383349179,1.047220310177444,1.013138620652269),
    @(1.03464234782153,1.028122220566489,1.037622906441491,1.026318811826916,1.048367288419864,1.013429720886154)
)

$bfCols = @(2,3,4,5,6)   # B,C,D,E,F
$inCols = @(9,10,11,12,13,14)  # I,J,K,L,M,N

for ($i = 0; $i -lt $bfData.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $bfData[$i]
    for ($j = 0; $j -lt $bfCols.Count; $j++) {
        $ws.Cells.Item($rowNum, $bfCols[$j]).Value = $rowVals[$j]
    }
}

for ($i = 0; $i -lt $inData.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $inData[$i]
    for ($j = 0; $j -lt $inCols.Count; $j++) {
        $ws.Cells.Item($rowNum, $inCols[$j]).Value = $rowVals[$j]
    }
}
